# Seeder workbook maintenance edit:
#  - Fix the misspelled "creteBy" header on the "content" sheet to "createBy".
#  - Leave the workbook with the "content" sheet active/selected (reflecting
#    the user's in-progress work on the content seeder), updating the
#    previously-active "category" sheet's lingering selection too.

$wb = $excel.ActiveWorkbook

# Fix the typo'd column header on the content sheet (H1: creteBy -> createBy)
$wsContent = $wb.Worksheets.Item("content")
$wsContent.Range("H1").Value = "createBy"

# The category sheet was the active tab before; update its remembered
# selection now that focus is moving away from it.
$wsCategory = $wb.Worksheets.Item("category")
$wsCategory.Range("F7").Select()

# Make the content sheet the active tab and set its selection.
$wsContent.Activate()
$wsContent.Range("H2").Select()
